# Generate Report for Handback
# Refresh the generated timestamps (and the zh-cn priority flag) that the
# handback report run produced on this pass.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# "Latest HO Xliff Generate Date" for the cd2c2209 / dfb8973f rows.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-27 04:16:34"
$wsOverview.Range("G5").Value = "2016-08-27 04:16:34"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority changed from "ht" (human translate) to "mt" (machine translate).
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
# Correspond Handoff Datetime
$wsZhCn.Range("H3").Value = "2016-08-27 04:16:30"
$wsZhCn.Range("H5").Value = "2016-08-27 04:16:30"
# Correspond Handback DateTime
$wsZhCn.Range("K3").Value = "2016-08-27 04:16:46"
$wsZhCn.Range("K5").Value = "2016-08-27 04:16:46"

# --- de-de sheet --------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
# Correspond Handback DateTime
$wsDeDe.Range("K3").Value = "2016-08-27 04:16:52"
$wsDeDe.Range("K5").Value = "2016-08-27 04:16:52"
